$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value looks like a plain number need to be forced
# to Text format first, otherwise Excel auto-converts the string to a number.
$numericPriceCells = @("D4", "D5", "D6", "D8", "D13", "D14", "D15", "D21", "D23", "D24", "D25", "D26", "D27", "D32", "D35", "D39", "D40", "D41", "D42", "D45", "D46", "D47", "D48", "D50")
foreach ($cellRef in $numericPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.335.27'
$ws.Range("E2").Value = '  +0.96%  '

$ws.Range("D3").Value = '2.446.94'
$ws.Range("E3").Value = '  +0.26%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").Value = '571.91'
$ws.Range("E5").Value = '  +0.90%  '

$ws.Range("D6").Value = '146.86'
$ws.Range("E6").Value = '  +0.84%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").Value = '0.540'
$ws.Range("E8").Value = '  +1.25%  '

$ws.Range("D9").Value = '2.442.23'
$ws.Range("E9").Value = '  -0.05%  '

$ws.Range("E10").Value = '  -0.22%  '

$ws.Range("E12").Value = '  -0.79%  '

$ws.Range("D13").Value = '0.356'
$ws.Range("E13").Value = '  +0.10%  '

$ws.Range("D14").Value = '27.05'
$ws.Range("E14").Value = '  +0.61%  '

$ws.Range("D15").Value = '0.0000179'
$ws.Range("E15").Value = '  -1.13%  '

$ws.Range("E16").Value = '  -0.20%  '

$ws.Range("D17").Value = '62.963.35'
$ws.Range("E17").Value = '  +0.54%  '

$ws.Range("D18").Value = '2.434.10'
$ws.Range("E18").Value = '  -0.10%  '

$ws.Range("E19").Value = '  +0.82%  '

$ws.Range("E20").Value = '  +5.50%  '

$ws.Range("D21").Value = '327.81'
$ws.Range("E21").Value = '  +1.19%  '

$ws.Range("E22").Value = '  +0.71%  '

$ws.Range("D23").Value = '2.10'
$ws.Range("E23").Value = '  +13.80%  '

$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.01%  '

$ws.Range("D25").Value = '65.38'

$ws.Range("D26").Value = '620.10'
$ws.Range("E26").Value = '  +6.22%  '

$ws.Range("D27").Value = '9.00'
$ws.Range("E27").Value = '  +4.50%  '

$ws.Range("E28").Value = '  +2.18%  '

$ws.Range("D29").Value = '2.561.18'
$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("E30").Value = '  +4.15%  '

$ws.Range("E31").Value = '  +0.38%  '

$ws.Range("D32").Value = '8.30'
$ws.Range("E32").Value = '  -1.40%  '

$ws.Range("E33").Value = '  -4.11%  '

$ws.Range("E34").Value = '  +0.36%  '

$ws.Range("D35").Value = '5.23'
$ws.Range("E35").Value = '  +8.02%  '

$ws.Range("E36").Value = '  +0.27%  '

$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("E38").Value = '  -0.65%  '

$ws.Range("D39").Value = '5.43'
$ws.Range("E39").Value = '  +0.75%  '

$ws.Range("D40").Value = '18.75'
$ws.Range("E40").Value = '  -0.28%  '

$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").Value = '2.70'
$ws.Range("E41").Value = '  +10.84%  '

$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").Value = '145.53'
$ws.Range("E42").Value = '  -1.76%  '

$ws.Range("E43").Value = '  -1.08%  '

$ws.Range("E44").Value = '  -0.19%  '

$ws.Range("D45").Value = '41.91'
$ws.Range("E45").Value = '  +0.65%  '

$ws.Range("D46").Value = '148.93'
$ws.Range("E46").Value = '  +0.16%  '

$ws.Range("D47").Value = '3.76'
$ws.Range("E47").Value = '  +2.33%  '

$ws.Range("D48").Value = '21.21'
$ws.Range("E48").Value = '  +3.44%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("D50").Value = '0.601'
$ws.Range("E50").Value = '  -0.07%  '

$ws.Range("E51").Value = '  +0.94%  '

# Restore default (Normal) style on the cells we temporarily reformatted,
# so the saved workbook keeps the original (unstyled) cell formatting.
foreach ($cellRef in $numericPriceCells) {
    $ws.Range($cellRef).Style = "Normal"
}